$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1338.75
$ws.Range("J121").Value = 1338.75
$ws.Range("L121").Value = 4016.25
$ws.Range("N121").Value = -7510.25

$ws.Range("H132").Value = 1072.0526
$ws.Range("I132").Value = 1132.0667
$ws.Range("J132").Value = 847
$ws.Range("K132").Value = 3396.2001
$ws.Range("L132").Value = 2541
$ws.Range("M132").Value = -866.2001
$ws.Range("N132").Value = -7601

$ws.Range("H137").Value = 2331.5454
$ws.Range("I137").Value = 1143.375
$ws.Range("K137").Value = 3430.125
$ws.Range("M137").Value = -880.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3038.1853
$ws.Range("I32").Value = 3039.6924
$ws.Range("K32").Value = 3039.6924
$ws.Range("M32").Value = -2752.6924

$ws.Range("H45").Value = 2959.4
$ws.Range("I45").Value = 2932.6667
$ws.Range("J45").Value = 2999.5
$ws.Range("K45").Value = 2932.6667
$ws.Range("L45").Value = 2999.5
$ws.Range("M45").Value = -2555.6667
$ws.Range("N45").Value = -3753.5

$ws.Range("H61").Value = 5998.143
$ws.Range("I61").Value = 6331.1665
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 6331.1665
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -6119.1665
$ws.Range("N61").Value = -4424

$ws.Range("H74").Value = 3564.1538
$ws.Range("I74").Value = 746.8
$ws.Range("J74").Value = 5325
$ws.Range("K74").Value = 746.8
$ws.Range("L74").Value = 5325
$ws.Range("M74").Value = 127.2
$ws.Range("N74").Value = -7073

$ws.Range("H77").Value = 3564.1538
$ws.Range("I77").Value = 746.8
$ws.Range("J77").Value = 5325
$ws.Range("K77").Value = 3734
$ws.Range("L77").Value = 26625
$ws.Range("M77").Value = 634
$ws.Range("N77").Value = -35361

$ws.Range("H122").Value = 3903.6667
$ws.Range("I122").Value = 3805.9092
$ws.Range("K122").Value = 11417.7276
$ws.Range("M122").Value = -8967.7276

$ws.Range("H132").Value = 3548.375
$ws.Range("I132").Value = 2897.8333
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 8693.499899999999
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -6163.499899999999
$ws.Range("N132").Value = -21560

$ws.Range("H135").Value = 149999.5
$ws.Range("J135").Value = 149999.5
$ws.Range("L135").Value = 149999.5
$ws.Range("N135").Value = -160139.5

$ws.Range("H136").Value = 5998.143
$ws.Range("I136").Value = 6331.1665
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 18993.4995
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -16443.4995
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 751.6
$ws.Range("I99").Value = 819.5
$ws.Range("K99").Value = 819.5
$ws.Range("M99").Value = 678.5

$ws.Range("H134").Value = 3143.6155
$ws.Range("I134").Value = 2897.0908
$ws.Range("K134").Value = 8691.2724
$ws.Range("M134").Value = -6156.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4278.8237
$ws.Range("I31").Value = 1288.3077
$ws.Range("J31").Value = 13998
$ws.Range("K31").Value = 1288.3077
$ws.Range("L31").Value = 13998
$ws.Range("M31").Value = -993.3077000000001
$ws.Range("N31").Value = -14588

$ws.Range("H34").Value = 4278.8237
$ws.Range("I34").Value = 1288.3077
$ws.Range("J34").Value = 13998
$ws.Range("K34").Value = 1288.3077
$ws.Range("L34").Value = 13998
$ws.Range("M34").Value = -1086.3077
$ws.Range("N34").Value = -14402

$ws.Range("H58").Value = 1562.8182
$ws.Range("I58").Value = 1532.4445
$ws.Range("J58").Value = 1699.5
$ws.Range("K58").Value = 1532.4445
$ws.Range("L58").Value = 1699.5
$ws.Range("M58").Value = -1329.4445
$ws.Range("N58").Value = -2105.5

$ws.Range("H62").Value = 3144.5
$ws.Range("I62").Value = 3144.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3144.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2520.5
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3144.5
$ws.Range("I65").Value = 3144.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15722.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -12602.5
$ws.Range("N65").ClearContents()

$ws.Range("H99").Value = 3501.8333
$ws.Range("I99").Value = 2669.2222
$ws.Range("J99").Value = 5999.6665
$ws.Range("K99").Value = 2669.2222
$ws.Range("L99").Value = 5999.6665
$ws.Range("M99").Value = -1171.2222
$ws.Range("N99").Value = -8995.666499999999

$ws.Range("H122").Value = 1002
$ws.Range("I122").Value = 1002
$ws.Range("K122").Value = 3006
$ws.Range("M122").Value = -556

$ws.Range("H126").Value = 3501.8333
$ws.Range("I126").Value = 2669.2222
$ws.Range("J126").Value = 5999.6665
$ws.Range("K126").Value = 8007.6666
$ws.Range("L126").Value = 17998.9995
$ws.Range("M126").Value = -5537.6666
$ws.Range("N126").Value = -22938.9995

$ws.Range("H136").Value = 1562.8182
$ws.Range("I136").Value = 1532.4445
$ws.Range("J136").Value = 1699.5
$ws.Range("K136").Value = 4597.333500000001
$ws.Range("L136").Value = 5098.5
$ws.Range("M136").Value = -2047.333500000001
$ws.Range("N136").Value = -10198.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6797.6
$ws.Range("I5").Value = 997.5
$ws.Range("K5").Value = 2992.5
$ws.Range("M5").Value = -2880.5

$ws.Range("H26").Value = 764.2143
$ws.Range("I26").Value = 633.2222
$ws.Range("J26").Value = 1000
$ws.Range("K26").Value = 1899.6666
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = -1611.6666
$ws.Range("N26").Value = -3576

$ws.Range("H60").Value = 910.5
$ws.Range("I60").Value = 849.25
$ws.Range("J60").Value = 971.75
$ws.Range("K60").Value = 2547.75
$ws.Range("L60").Value = 2915.25
$ws.Range("M60").Value = -2296.75
$ws.Range("N60").Value = -3417.25

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H135").Value = 6797.6
$ws.Range("I135").Value = 997.5
$ws.Range("K135").Value = 8977.5
$ws.Range("M135").Value = -6442.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4575
$ws.Range("I102").Value = 4575
$ws.Range("K102").Value = 4575
$ws.Range("M102").Value = -2953

$ws.Range("H126").Value = 8797.200000000001
$ws.Range("I126").Value = 7996.3335
$ws.Range("K126").Value = 23989.0005
$ws.Range("M126").Value = -21519.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 18750.5
$ws.Range("J6").Value = 19999
$ws.Range("L6").Value = 19999
$ws.Range("N6").Value = -20229

$ws.Range("H132").Value = 3478.054
$ws.Range("I132").Value = 2829.75
$ws.Range("J132").Value = 4240.7646
$ws.Range("K132").Value = 8489.25
$ws.Range("L132").Value = 12722.2938
$ws.Range("M132").Value = -5959.25
$ws.Range("N132").Value = -17782.2938

$ws.Range("H136").Value = 11320.3125
$ws.Range("I136").Value = 11429.2
$ws.Range("J136").Value = 9687
$ws.Range("K136").Value = 34287.60000000001
$ws.Range("L136").Value = 29061
$ws.Range("M136").Value = -31737.60000000001
$ws.Range("N136").Value = -34161
